# Daily attendance processing - 2026-01-21 08:07:31
# Applies the day's attendance-recording pass to the "Session Analysis Results" sheet:
#   - Session 23 (26/10/2025) for groups B2A/B2B/B2C gets recorded (B2A was fully
#     "Not Recorded"; B2B and B2C get their headcounts corrected).
#   - "Recorded By" lists get the author's "System" entry reordered to the front.
#   - Class Statistics / per-group summary numbers refresh accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Recorded-By cell text: move "System" to the front of the comma list while
#    keeping the relative order of the remaining names (purely a re-ordering).
# ---------------------------------------------------------------------------
$recordedByUpdates = @{
    "G2"   = "System, system, backup@backdoor.com"
    "G3"   = "System, dnasr281@gmail.com"
    "G5"   = "System, backup@backdoor.com"
    "G6"   = "System, dnasr281@gmail.com"
    "G8"   = "System, backup@backdoor.com"
    "G28"  = "System, system, backup@backdoor.com"
    "G29"  = "System, dnasr281@gmail.com"
    "G31"  = "System, backup@backdoor.com"
    "G32"  = "System, dnasr281@gmail.com"
    "G34"  = "System, backup@backdoor.com"
    "G54"  = "System, system, backup@backdoor.com"
    "G55"  = "System, dnasr281@gmail.com"
    "G57"  = "System, backup@backdoor.com"
    "G58"  = "System, dnasr281@gmail.com"
    "G60"  = "System, backup@backdoor.com"
    "G80"  = "System, backup@backdoor.com"
    "G81"  = "System, backup@backdoor.com"
    "G82"  = "System, backup@backdoor.com"
    "G106" = "System, backup@backdoor.com"
    "G107" = "System, backup@backdoor.com"
    "G108" = "System, backup@backdoor.com"
    "G132" = "System, backup@backdoor.com"
    "G133" = "System, backup@backdoor.com"
    "G134" = "System, backup@backdoor.com"
}
foreach ($addr in $recordedByUpdates.Keys) {
    $ws.Range($addr).Value = $recordedByUpdates[$addr]
}

# ---------------------------------------------------------------------------
# 2) Row 24 (Year4 / B2A / session 23) moves from "Not Recorded" (special
#    highlighted style) to a normally-styled, fully recorded row. Re-apply the
#    ordinary data-row formatting (copied from row 23, which already carries
#    it) and then refresh the text that actually changed.
# ---------------------------------------------------------------------------
$ws.Range("A23:I23").Copy()
$ws.Range("A24:I24").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("G24").Value = "System, dnasr281@gmail.com"
$ws.Range("H24").Value = "27/52"
$ws.Range("I24").Value = "Recorded"

# ---------------------------------------------------------------------------
# 3) Matching session-23 headcounts for the sibling groups (B2B row 50, B2C
#    row 76) also update now that the session has been processed.
# ---------------------------------------------------------------------------
$ws.Range("H50").Value = "41/57"
$ws.Range("H76").Value = "35/55"

# ---------------------------------------------------------------------------
# 4) Class Statistics block (K/L columns near the top) refreshes.
# ---------------------------------------------------------------------------
$ws.Range("L6").Value = 156   # Recorded Sessions
$ws.Range("L7").Value = 0     # Missing Sessions

# Percentage cells are stored as literal text (e.g. "100.0%"), not numbers -
# a plain `.Value = "100.0%"` assignment gets auto-coerced by Excel into a
# numeric percentage (and a new NumberFormat style). Instead, enter it as a
# text-returning formula, then collapse the formula to its literal value with
# a values-only paste: the cell keeps its original style untouched and ends
# up holding the plain text, exactly like the source data.
function Set-TextValue($addr, $text) {
    $escaped = $text.Replace('"', '""')
    $ws.Range($addr).Formula = '="' + $escaped + '"'
    $ws.Range($addr).Copy()
    $ws.Range($addr).PasteSpecial(-4163)
    $excel.CutCopyMode = $false
}

Set-TextValue "L9"  "100.0%"   # Coverage %
Set-TextValue "L10" "67.6%"    # Average Attendance %

# ---------------------------------------------------------------------------
# 5) Per-group summary table (rows 15-17, columns M-S) refreshes for the B2A
#    group row (row 15) and the Average Attendance % for B2B/B2C (rows 16/17).
# ---------------------------------------------------------------------------
$ws.Range("O15").Value = 26
$ws.Range("P15").Value = 0

Set-TextValue "R15" "100.0%"
Set-TextValue "S15" "68.0%"
Set-TextValue "S16" "64.6%"
Set-TextValue "S17" "61.3%"

# ---------------------------------------------------------------------------
# 6) Column I narrows slightly (14 -> 10 characters of raw stored width).
#    ColumnWidth is expressed in "characters"; the engine stores raw width as
#    ColumnWidth + 5/6, so subtract that offset to land exactly on width=10.
# ---------------------------------------------------------------------------
$ws.Columns.Item(9).ColumnWidth = 10 - 5/6
